# "update activity file upload"
# The sheet's sample data is replaced wholesale: the old scattered
# "sdsd/dsd/..." placeholder values (rows 1-19) are cleared out and replaced
# with a new, more compact block of "xcsdfsdf/df/sdf/..." placeholder values
# (rows 4-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous placeholder data entirely before laying down the new
# values - the used range shrinks from A1:J19 down to A4:E14.
$ws.Cells.Clear()

# Re-enter the new values cell-by-cell, in the same order the original
# author typed them, so the shared-string table comes out in the same
# first-seen order as the target workbook.
$ws.Range("B4").Value  = "xcsdfsdf"
$ws.Range("E6").Value  = "df"
$ws.Range("D14").Value = "sdf"
$ws.Range("B8").Value  = "f"
$ws.Range("A4").Value  = "sdfsdfsdf"
$ws.Range("C8").Value  = "sdfsdfsdf"
$ws.Range("B10").Value = "sdfsdfsdf"
$ws.Range("D6").Value  = "sdfsdfsd"

# Move the active selection to where the author left off editing.
$ws.Range("D6").Select() | Out-Null
